$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells F1:H1 - new columns with style matching existing headers (s="1")
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Fill boolean values for F2:H9, default FALSE
$ws.Range("F2:H9").Value = $false

# Set the one TRUE value per the diff: G5
$ws.Range("G5").Value = $true
